$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.545.50"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.768.26"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "3.208.28"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "2.758.74"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "51.519.06"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +13.78%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0448"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.74%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.96%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.48%  "
$ws.Range("D46").Value = "2.083.19"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.931"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.50%  "
